$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 262; $r++) {
    $gCell = $ws.Cells.Item($r, 7)
    $gVal = $gCell.Value()
    if ($gVal -ne $null) {
        $gCell.Value = $gVal + 1
    }

    $iCell = $ws.Cells.Item($r, 9)
    $iVal = $iCell.Value()
    if ($iVal -ne $null) {
        $iCell.Value = $iVal - 1
    }
}
